$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing content so stale shared strings are dropped
$ws.Cells.Clear()

# Establish shared-string first-use order: Input, Output, GPT-4o, GPT-4o Batched, Token
$ws.Range("B1").Value = "Input"
$ws.Range("F1").Value = "Output"
$ws.Range("A1").Value = "GPT-4o"
$ws.Range("A2").Value = "GPT-4o Batched"
$ws.Range("E1").Value = "Token"

# Row 1 - GPT-4o
$ws.Range("C1").Value = 5
$ws.Range("D1").Value = 1000000
$ws.Range("G1").Value = 15
$ws.Range("H1").Value = 1000000
$ws.Range("I1").Value = "Token"

# Row 2 - GPT-4o Batched
$ws.Range("B2").Value = "Input"
$ws.Range("C2").Value = 2.5
$ws.Range("D2").Value = 1000000
$ws.Range("E2").Value = "Token"
$ws.Range("F2").Value = "Output"
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 1000000
$ws.Range("I2").Value = "Token"

# Move the active selection like the author left it
$ws.Range("F3").Select()
